$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

$data = @(
    @("options", "OPTIONS"),
    @("music", "MUSIC"),
    @("sound", "SOUND"),
    @("speech", "SPEECH"),
    @("on", "ON"),
    @("off", "OFF"),
    @("close", "CLOSE"),
    @("newGame", "NEW GAME"),
    @("continue", "CONTINUE")
)

$row = 4
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

$ws.Range("B12").Select()
